$d = $word.ActiveDocument

# 1. Update the title paragraph text (append a trailing space)
$d.Paragraphs.Item(1).Range.Text = "BEHAVIOURAL "

# 2. Remove paragraphs 2 through 9 (the diagram shapes block) entirely, then
#    re-insert the replacement two-paragraph block (Vehicle/Security box row,
#    Access/Password box row + connectors) via InsertXML so the VML <v:pict>
#    shapes come along with the paragraph markup.
$startPos = $d.Paragraphs.Item(2).Range.Start
$endPos = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$d.Range($startPos, $endPos).Delete()

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:rect id="_x0000_s1033" style="position:absolute;margin-left:213.75pt;margin-top:18.05pt;width:105pt;height:24pt;z-index:251664384"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Vehicle</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:rect id="_x0000_s1026" style="position:absolute;margin-left:30.75pt;margin-top:21.8pt;width:106.5pt;height:24pt;z-index:251658240"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Security</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shapetype id="_x0000_t32" coordsize="21600,21600" o:spt="32" o:oned="t" path="m,l21600,21600e" filled="f"><v:path arrowok="t" fillok="f" o:connecttype="none"/><o:lock v:ext="edit" shapetype="t"/></v:shapetype><v:shape id="_x0000_s1032" type="#_x0000_t32" style="position:absolute;margin-left:267pt;margin-top:16.6pt;width:0;height:51.75pt;flip:y;z-index:251663360" o:connectortype="straight"><v:stroke endarrow="block"/></v:shape></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:rect id="_x0000_s1028" style="position:absolute;margin-left:215.25pt;margin-top:68.35pt;width:103.5pt;height:27pt;z-index:251660288"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Access</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shape id="_x0000_s1030" type="#_x0000_t32" style="position:absolute;margin-left:137.25pt;margin-top:81.85pt;width:78pt;height:.75pt;z-index:251662336" o:connectortype="straight"><v:stroke endarrow="block"/></v:shape></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:shape id="_x0000_s1029" type="#_x0000_t32" style="position:absolute;margin-left:81.75pt;margin-top:20.35pt;width:0;height:51.75pt;z-index:251661312" o:connectortype="straight"><v:stroke endarrow="block"/></v:shape></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict><v:rect id="_x0000_s1027" style="position:absolute;margin-left:30.75pt;margin-top:72.1pt;width:106.5pt;height:23.25pt;z-index:251659264"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Password</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></w:r></w:p>'

$d.Paragraphs.Item(2).Range.InsertXML($xml)
